$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 24 (the ERROR_TEST row),
# shifting it (and everything below) down by one.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new "ERROR_LOGS_NOT_FOUND" error code.
$ws.Range("A24").Value = 1022
$ws.Range("B24").Value = "ERROR_LOGS_NOT_FOUND"
$ws.Range("C24").Value = "general"
$ws.Range("D24").Value = "日志读取失败"

# Match the selection recorded in the saved workbook.
$ws.Range("B24").Select()
